# Insert a new data row at row 27 (pushes existing rows 27-53 down to 28-54)
# and populate it with a new weekly price record, matching the rest of the
# "Hortaliza, Femacal de La Calera - Arveja Verde" dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(27).Insert()

$ws.Cells.Item(27, 1).Value = 3
$ws.Cells.Item(27, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(27, 3).Value = "Coquimbo"
$ws.Cells.Item(27, 4).Value = "2022-02-09"
$ws.Cells.Item(27, 5).Value = 5
$ws.Cells.Item(27, 6).Value = 100112022
$ws.Cells.Item(27, 7).Value = "Arveja Verde"
$ws.Cells.Item(27, 8).Value = "Perfection"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 50
$ws.Cells.Item(27, 11).Value = 23000
$ws.Cells.Item(27, 12).Value = 24000
$ws.Cells.Item(27, 13).Value = 23600
$ws.Cells.Item(27, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(27, 15).Value = "Región Metropolitana"
$ws.Cells.Item(27, 16).Value = 944
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"
